$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

$ws.Range("B2").Value = -4.889486674598956
$ws.Range("C2").Value = 2.5042050373577567
$ws.Range("D2").Value = 0.67431139478685509
$ws.Range("E2").Value = -0.014437625172774915

$ws.Range("B3").Value = 2.8609691566184168
$ws.Range("C3").Value = 3.7885273178131698
$ws.Range("D3").Value = 5.9981617106704093
$ws.Range("E3").Value = -3.2221819687262965

$ws.Range("B1:E3").Select()
